{"js": "// Update the worksheet date and the 25 multiplication problems to the\n// new \"output generated at c986bee\" values.\nconst replacements = [\n  [\"2024-10-26 Saturday\", \"2024-10-27 Sunday\"],\n  [\"221\u00d78=\", \"401\u00d77=\"],\n  [\"455\u00d72=\", \"552\u00d77=\"],\n  [\"295\u00d72=\", \"192\u00d77=\"],\n  [\"818\u00d77=\", \"568\u00d75=\"],\n  [\"350\u00d74=\", \"348\u00d73=\"],\n  [\"969\u00d73=\", \"564\u00d78=\"],\n  [\"979\u00d79=\", \"898\u00d74=\"],\n  [\"287\u00d77=\", \"627\u00d74=\"],\n  [\"770\u00d73=\", \"290\u00d76=\"],\n  [\"195\u00d77=\", \"567\u00d79=\"],\n  [\"608\u00d74=\", \"540\u00d73=\"],\n  [\"649\u00d75=\", \"298\u00d79=\"],\n  [\"207\u00d79=\", \"674\u00d74=\"],\n  [\"924\u00d79=\", \"762\u00d79=\"],\n  [\"624\u00d79=\", \"279\u00d73=\"],\n  [\"876\u00d75=\", \"705\u00d72=\"],\n  [\"720\u00d77=\", \"606\u00d79=\"],\n  [\"881\u00d75=\", \"578\u00d73=\"],\n  [\"161\u00d79=\", \"387\u00d79=\"],\n  [\"830\u00d75=\", \"294\u00d74=\"],\n  [\"354\u00d72=\", \"637\u00d79=\"],\n  [\"321\u00d78=\", \"206\u00d74=\"],\n  [\"741\u00d72=\", \"251\u00d73=\"],\n  [\"437\u00d73=\", \"595\u00d75=\"],\n  [\"349\u00d76=\", \"558\u00d72=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update the worksheet date and the 25 multiplication problems to the\n# new \"output generated at c986bee\" values.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2024-10-26 Saturday\", \"2024-10-27 Sunday\"),\n    @(\"221\u00d78=\", \"401\u00d77=\"),\n    @(\"455\u00d72=\", \"552\u00d77=\"),\n    @(\"295\u00d72=\", \"192\u00d77=\"),\n    @(\"818\u00d77=\", \"568\u00d75=\"),\n    @(\"350\u00d74=\", \"348\u00d73=\"),\n    @(\"969\u00d73=\", \"564\u00d78=\"),\n    @(\"979\u00d79=\", \"898\u00d74=\"),\n    @(\"287\u00d77=\", \"627\u00d74=\"),\n    @(\"770\u00d73=\", \"290\u00d76=\"),\n    @(\"195\u00d77=\", \"567\u00d79=\"),\n    @(\"608\u00d74=\", \"540\u00d73=\"),\n    @(\"649\u00d75=\", \"298\u00d79=\"),\n    @(\"207\u00d79=\", \"674\u00d74=\"),\n    @(\"924\u00d79=\", \"762\u00d79=\"),\n    @(\"624\u00d79=\", \"279\u00d73=\"),\n    @(\"876\u00d75=\", \"705\u00d72=\"),\n    @(\"720\u00d77=\", \"606\u00d79=\"),\n    @(\"881\u00d75=\", \"578\u00d73=\"),\n    @(\"161\u00d79=\", \"387\u00d79=\"),\n    @(\"830\u00d75=\", \"294\u00d74=\"),\n    @(\"354\u00d72=\", \"637\u00d79=\"),\n    @(\"321\u00d78=\", \"206\u00d74=\"),\n    @(\"741\u00d72=\", \"251\u00d73=\"),\n    @(\"437\u00d73=\", \"595\u00d75=\"),\n    @(\"349\u00d76=\", \"558\u00d72=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute($find.Text, $true, $true, $false, $false, $false, $true, \"wdFindContinue\", $false, $find.Replacement.Text, \"wdReplaceAll\")\n}\n"}
